$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column (04-dec) before the
#     01-oct. block (currently column EG), shifting all following
#     columns one to the right (through former FK -> FL). ---
$ws1 = $wb.Worksheets.Item("Prix Spot")
$ws1.Range("EG1:EG25").EntireColumn.Insert()
$ws1.Range("EG1").Value = "04-dec"
$ws1.Range("EG2:EG25").Value = "-"

# --- Sheet "Gaz": append new row 167 for 2025-12-02. ---
$ws2 = $wb.Worksheets.Item("Gaz")
$ws2.Range("A167").Value = "'2025-12-02"
$ws2.Range("A167").ClearFormats()
$ws2.Range("B167").Value = 26.895

# --- Sheet "CO2": append new row 167 for 2025-12-02. ---
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Range("A167").Value = "'2025-12-02"
$ws3.Range("A167").ClearFormats()
$ws3.Range("B167").Value = 81.65
